# This script appends the "Статья ... Скалько ..." paragraph plus the
# trailing signature block to the recommendation letter, adding Skalko as
# a co-author of the reviewed article (per the commit message).

$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph (right after
# the "Рекомендательное письмо" heading). That paragraph's own closing
# paragraph mark must remain the very last mark in the document body, so
# we insert the seven new paragraphs' worth of OOXML in front of it and
# then fold that trailing empty paragraph into the new final paragraph.
$tailPara = $d.Paragraphs.Last
$tailRange = $tailPara.Range

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:contextualSpacing/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:contextualSpacing/>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Статья А.В. Шевченко</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> И.В. </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Цыбулина</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> и Ю.И. Скалько</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> «</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Оптимизационный алгоритм решения задачи о фазовом</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>равновесии</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">» была рассмотрена на научном семинаре лаборатории </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>флюидодинамики</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> и </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>сейсмоакустики</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> МФТИ и рекомендована для публикации в журнале «Труды МФТИ».</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:contextualSpacing/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:contextualSpacing/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:contextualSpacing/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:contextualSpacing/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="708"/>
        <w:contextualSpacing/>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Зав. Лабораторией</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:tab/>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Колдоба</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> А.В.</w:t>
      </w:r>

    </w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$tailRange.InsertXML($xml)

# After the insert, the document looks like:
#   ... "Рекомендательное письмо" ...
#   <new paragraph 1 - empty, justified, first-line indent>
#   <new paragraph 2 - "Статья А.В. Шевченко, ... Скалько ...">
#   <new paragraph 3..6 - empty, first-line indent>
#   <new paragraph 7 - "Зав. Лабораторией ... Колдоба А.В.">
#   <original trailing empty paragraph>   <-- still last, unchanged
#
# Fold the original trailing empty paragraph into the freshly inserted
# "Зав. Лабораторией" paragraph by deleting the paragraph mark that
# separates them; Word keeps the *later* paragraph's pPr when a mark is
# deleted this way, which is exactly the (contextualSpacing + rPr sz28)
# pPr the final paragraph needs.
$signaturePara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$markRange = $d.Range($signaturePara.Range.End - 1, $signaturePara.Range.End)
$markRange.Delete()

# Finally, give the merged final paragraph its first-line indent (708
# twips = 35.4 pt), matching the other new paragraphs.
$finalPara = $d.Paragraphs.Last
$finalPara.Range.ParagraphFormat.FirstLineIndent = 708 / 20.0
